# "Generate Report for Handback" - refresh the handback timestamps
# recorded in the report workbook.
#
# Overview!G2 ("Latest HO Xliff Generate Date" for the first row) and
# de-de!H2 ("Correspond Handoff Datetime" for the first row) both held
# the same shared string ("2016-08-19 23:08:26"); it becomes
# "2016-08-19 23:09:11" for both.
#
# zh-cn!H2 / zh-cn!K2 (Correspond Handoff/Handback Datetime for the
# first row) move from 23:08:20 / 23:08:37 to 23:09:06 / 23:09:24.
#
# de-de!K2 (Correspond Handback Datetime for the first row) moves from
# 23:08:43 to 23:09:31.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-19 23:09:11"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-19 23:09:06"
$wsZhCn.Range("K2").Value = "2016-08-19 23:09:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-19 23:09:11"
$wsDeDe.Range("K2").Value = "2016-08-19 23:09:31"
